# Auto-generated edit script applying the committed numeric updates
# to the Anima_Profits workbook (multi-sheet cell value corrections).
$wb = $excel.ActiveWorkbook

# --- ALC!row 92 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 22223294
$ws.Range("I92").Value = 27778218
$ws.Range("J92").Value = 3600
$ws.Range("K92").Value = 27778218
$ws.Range("L92").Value = 3600
$ws.Range("M92").Value = -27776970
$ws.Range("N92").Value = -6096

# --- ALC!row 138 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2309.49
$ws.Range("I138").Value = 1373.5
$ws.Range("J138").Value = 2348.4895
$ws.Range("K138").Value = 4120.5
$ws.Range("L138").Value = 7045.468500000001
$ws.Range("M138").Value = 1019.5
$ws.Range("N138").Value = -17325.4685

# --- ARM!row 28 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 3414.2
$ws.Range("I28").Value = 3414.2
$ws.Range("K28").Value = 3414.2
$ws.Range("M28").Value = -3222.2

# --- ARM!row 32 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1123243.4
$ws.Range("I32").Value = 1289750.6
$ws.Range("J32").Value = 40946.5
$ws.Range("K32").Value = 1289750.6
$ws.Range("L32").Value = 40946.5
$ws.Range("M32").Value = -1289463.6
$ws.Range("N32").Value = -41520.5

# --- ARM!row 99 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H99").Value = 3414.2
$ws.Range("I99").Value = 3414.2
$ws.Range("K99").Value = 3414.2
$ws.Range("M99").Value = -419.1999999999998

# --- ARM!row 133 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 56000
$ws.Range("J133").Value = 56000
$ws.Range("L133").Value = 56000
$ws.Range("N133").Value = -61060

# --- CRP!row 16 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1880
$ws.Range("I16").Value = 1250
$ws.Range("J16").Value = 2510
$ws.Range("K16").Value = 1250
$ws.Range("L16").Value = 2510
$ws.Range("M16").Value = -963
$ws.Range("N16").Value = -3084

# --- CRP!row 22 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 521.125
$ws.Range("I22").Value = 494.66666
$ws.Range("J22").Value = 600.5
$ws.Range("K22").Value = 494.66666
$ws.Range("L22").Value = 600.5
$ws.Range("M22").Value = -144.66666
$ws.Range("N22").Value = -1300.5

# --- CRP!row 113 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1880
$ws.Range("I113").Value = 1250
$ws.Range("J113").Value = 2510
$ws.Range("K113").Value = 1250
$ws.Range("L113").Value = 2510
$ws.Range("M113").Value = 920
$ws.Range("N113").Value = -6850

# --- CRP!row 122 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1346.091
$ws.Range("I122").Value = 1346.091
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4038.273
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1588.273
$ws.Range("N122").ClearContents()

# --- CUL!row 44 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 481.625
$ws.Range("I44").Value = 387.5
$ws.Range("J44").Value = 575.75
$ws.Range("K44").Value = 1162.5
$ws.Range("L44").Value = 1727.25
$ws.Range("M44").Value = -764.5
$ws.Range("N44").Value = -2523.25

# --- CUL!row 107 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 446.47058
$ws.Range("I107").Value = 236.66667
$ws.Range("J107").Value = 950
$ws.Range("K107").Value = 710.00001
$ws.Range("L107").Value = 2850
$ws.Range("M107").Value = 1209.99999
$ws.Range("N107").Value = -6690

# --- CUL!row 134 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 6719.8184
$ws.Range("I134").Value = 3605.4443
$ws.Range("J134").Value = 8875.923000000001
$ws.Range("K134").Value = 10816.3329
$ws.Range("L134").Value = 26627.769
$ws.Range("M134").Value = -5746.332900000001
$ws.Range("N134").Value = -36767.769

# --- GSM!row 70 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5026.646
$ws.Range("I70").Value = 5071.0513
$ws.Range("J70").Value = 4834.222
$ws.Range("K70").Value = 5071.0513
$ws.Range("L70").Value = 4834.222
$ws.Range("M70").Value = -4801.0513
$ws.Range("N70").Value = -5374.222

# --- GSM!row 73 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5026.646
$ws.Range("I73").Value = 5071.0513
$ws.Range("J73").Value = 4834.222
$ws.Range("K73").Value = 5071.0513
$ws.Range("L73").Value = 4834.222
$ws.Range("M73").Value = -4135.0513
$ws.Range("N73").Value = -6706.222

# --- GSM!row 113 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 93650.45
$ws.Range("I113").Value = 113795
$ws.Range("K113").Value = 113795
$ws.Range("M113").Value = -111625

# --- GSM!row 122 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2500
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5050
$ws.Range("N122").ClearContents()

# --- GSM!row 132 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3836.0625
$ws.Range("I132").Value = 3039.111
$ws.Range("J132").Value = 4860.7144
$ws.Range("K132").Value = 9117.332999999999
$ws.Range("L132").Value = 14582.1432
$ws.Range("M132").Value = -6587.332999999999
$ws.Range("N132").Value = -19642.1432

# --- GSM!row 136 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 29865.2
$ws.Range("J136").Value = 20108.666
$ws.Range("L136").Value = 60325.99800000001
$ws.Range("N136").Value = -65425.99800000001

# --- LTW!row 40 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3055.1875
$ws.Range("I40").Value = 3052.9092
$ws.Range("J40").Value = 3060.2
$ws.Range("K40").Value = 3052.9092
$ws.Range("L40").Value = 3060.2
$ws.Range("M40").Value = -2916.9092
$ws.Range("N40").Value = -3332.2

# --- LTW!row 82 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 50003500
$ws.Range("I82").Value = 100002600
$ws.Range("J82").Value = 4400.4
$ws.Range("K82").Value = 100002600
$ws.Range("L82").Value = 4400.4
$ws.Range("M82").Value = -100002239
$ws.Range("N82").Value = -5122.4

# --- LTW!row 85 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 50003500
$ws.Range("I85").Value = 100002600
$ws.Range("J85").Value = 4400.4
$ws.Range("K85").Value = 100002600
$ws.Range("L85").Value = 4400.4
$ws.Range("M85").Value = -100001352
$ws.Range("N85").Value = -6896.4

# --- LTW!row 122 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2938.875
$ws.Range("I122").Value = 2902.2
$ws.Range("K122").Value = 8706.599999999999
$ws.Range("M122").Value = -6256.599999999999

# --- WVR!row 107 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 913.90625
$ws.Range("I107").Value = 808.65
$ws.Range("K107").Value = 2425.95
$ws.Range("M107").Value = -505.9499999999998

# --- WVR!row 122 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1943.7778
$ws.Range("I122").Value = 1903.1538
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 5709.4614
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -3259.4614
$ws.Range("N122").Value = -13900

# --- WVR!row 126 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 677.7778
$ws.Range("I126").Value = 677.7778
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2033.3334
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = 436.6666
$ws.Range("N126").ClearContents()

# --- WVR!row 132 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6582569
$ws.Range("I132").Value = 4427.636
$ws.Range("J132").Value = 15627513
$ws.Range("K132").Value = 13282.908
$ws.Range("L132").Value = 46882539
$ws.Range("M132").Value = -10752.908
$ws.Range("N132").Value = -46887599

# --- WVR!row 136 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6311.6665
$ws.Range("I136").Value = 7797.778
$ws.Range("J136").Value = 5197.0835
$ws.Range("K136").Value = 23393.334
$ws.Range("L136").Value = 15591.2505
$ws.Range("M136").Value = -20843.334
$ws.Range("N136").Value = -20691.2505

# --- WVR!row 141 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 63365
$ws.Range("J141").Value = 63365
$ws.Range("L141").Value = 63365
$ws.Range("N141").Value = -73725
